$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the match data between row 61 and row 62 (columns F:V only; ---
# --- columns A:E - Indice/pais/torneio/temporada/data_partida - stay put) ---
$row61 = @($ws.Range("F61").Value2, $ws.Range("G61").Value2, $ws.Range("H61").Value2, $ws.Range("I61").Value2, $ws.Range("J61").Value2, $ws.Range("K61").Value2, $ws.Range("L61").Value2, $ws.Range("M61").Value2, $ws.Range("N61").Value2, $ws.Range("O61").Value2, $ws.Range("P61").Value2, $ws.Range("Q61").Value2, $ws.Range("R61").Value2, $ws.Range("S61").Value2, $ws.Range("T61").Value2, $ws.Range("U61").Value2, $ws.Range("V61").Value2)
$row62 = @($ws.Range("F62").Value2, $ws.Range("G62").Value2, $ws.Range("H62").Value2, $ws.Range("I62").Value2, $ws.Range("J62").Value2, $ws.Range("K62").Value2, $ws.Range("L62").Value2, $ws.Range("M62").Value2, $ws.Range("N62").Value2, $ws.Range("O62").Value2, $ws.Range("P62").Value2, $ws.Range("Q62").Value2, $ws.Range("R62").Value2, $ws.Range("S62").Value2, $ws.Range("T62").Value2, $ws.Range("U62").Value2, $ws.Range("V62").Value2)

$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "61").Value = $row62[$i]
    $ws.Range($cols[$i] + "62").Value = $row61[$i]
}

# --- Append two new match rows (123, 124) copying row 122's formatting ---
$ws.Range("A122:V122").Copy()
$ws.Range("A123:V123").PasteSpecial(-4122)
$ws.Range("A122:V122").Copy()
$ws.Range("A124:V124").PasteSpecial(-4122)

# Row 123: Crotone 2 - 1 Monterosi
$ws.Range("A123").Value = 122
$ws.Range("B123").Value = "italy"
$ws.Range("C123").Value = "serie-c-group-c"
$ws.Range("D123").Value = "2023-2024"
$ws.Range("E123").Value = 45242.58333333334
$ws.Range("F123").Value = "Crotone"
$ws.Range("G123").Value = 2
$ws.Range("H123").Value = "Monterosi"
$ws.Range("I123").Value = 1
$ws.Range("J123").Value = 1.33
$ws.Range("K123").Value = "09/11/2023 09:13"
$ws.Range("L123").Value = 1.34
$ws.Range("M123").Value = "09/11/2023 14:58"
$ws.Range("N123").Value = 4.71
$ws.Range("O123").Value = "09/11/2023 09:13"
$ws.Range("P123").Value = 5.03
$ws.Range("Q123").Value = "12/11/2023 12:59"
$ws.Range("R123").Value = 8.039999999999999
$ws.Range("S123").Value = "09/11/2023 09:13"
$ws.Range("T123").Value = 9.02
$ws.Range("U123").Value = "12/11/2023 12:15"
$ws.Range("V123").Value = "https://www.betexplorer.com/football/italy/serie-c-group-c/crotone-monterosi/zs36Jg8e/"

# Row 124: Turris 0 - 1 Monopoli
$ws.Range("A124").Value = 123
$ws.Range("B124").Value = "italy"
$ws.Range("C124").Value = "serie-c-group-c"
$ws.Range("D124").Value = "2023-2024"
$ws.Range("E124").Value = 45242.58333333334
$ws.Range("F124").Value = "Turris"
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = "Monopoli"
$ws.Range("I124").Value = 1
$ws.Range("J124").Value = 2.16
$ws.Range("K124").Value = "09/11/2023 09:13"
$ws.Range("L124").Value = 2.33
$ws.Range("M124").Value = "12/11/2023 13:46"
$ws.Range("N124").Value = 2.9
$ws.Range("O124").Value = "09/11/2023 09:13"
$ws.Range("P124").Value = 3.35
$ws.Range("Q124").Value = "12/11/2023 13:46"
$ws.Range("R124").Value = 3.41
$ws.Range("S124").Value = "09/11/2023 09:13"
$ws.Range("T124").Value = 3.03
$ws.Range("U124").Value = "12/11/2023 13:42"
$ws.Range("V124").Value = "https://www.betexplorer.com/football/italy/serie-c-group-c/turris-monopoli/l6hiPqdE/"
